# "Generate Report for Handoff"
#
# A second file (ec26f46a-580f-4437-9711-2f0e87f9a527.md) has gone through the
# handoff pipeline alongside the existing d8f2571c-... file, whose own guid
# has rolled forward to d9d4f421-203f-4577-85e6-53a0a07e7898. Every sheet
# gains a row for the new file and keeps the ".localization-config" /
# "not to be localized" row last.

$wb = $excel.ActiveWorkbook

# d8f2571c-...  (old source guid)      -> d9d4f421-... (rolled forward)
# ec26f46a-...                          -> brand-new source file, appended
$newGuid = "d9d4f421-203f-4577-85e6-53a0a07e7898"
$addedGuid = "ec26f46a-580f-4437-9711-2f0e87f9a527"

# content hash embedded in the handoff/handback .xlf file names; shared by
# both the zh-cn and de-de variants of a given source file
$newHash = "260536dd4e1c9f7f8b9a0f36b95538790ad3e2b1"
$addedHash = "cbf522627f973f49b2d45468156752f713f79521"

$zhTime = "2016-02-22 05:02:57"
$deTime = "2016-02-22 05:03:10"

$epoch = "0001-01-01 00:00:00"

function Set-MdHyperlink($ws, $cellRef, $guid, $repoBase) {
    $target = "$repoBase/e2e/$guid.md"
    $ws.Range($cellRef).Value = "$guid.md"
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", "$guid.md") | Out-Null
}

function Set-XlfHyperlink($ws, $cellRef, $guid, $hash, $lang, $handoffBase) {
    $fileName = "$guid.$hash.$lang.xlf"
    $target = "$handoffBase/$fileName"
    $ws.Range($cellRef).Value = $fileName
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $fileName) | Out-Null
}

$overviewRepo = "https://github.com/OpenLocalizationTest/oltest/blob/1b58219126d4f5cb5b4679eca9d8a2c9290108b7"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9aab5513670a8d36b2fe9c1aa2bf7b2c462ab656/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9ba647320b2a68a60e25c0ae90d83d909627390/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Hyperlinks.Delete()

# push the config row down from 3 -> 4, then fill in the new row 3
$wsOv.Range("A4").Value = ".localization-config"
$wsOv.Range("B4").Value = "Not to be localized"
$wsOv.Range("C4").Value = "Not to be localized"

$wsOv.Range("B3").Value = "Ready for handoff"
$wsOv.Range("C3").Value = "Ready for handoff"

Set-MdHyperlink $wsOv "A2" $newGuid $overviewRepo
Set-MdHyperlink $wsOv "A3" $addedGuid $overviewRepo
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "$overviewRepo/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A4").Value = ".localization-config"
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $zhTime
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("D2").Value = $zhTime

Set-MdHyperlink $wsZh "A2" $newGuid $overviewRepo
Set-XlfHyperlink $wsZh "C2" $newGuid $newHash "zh-cn" $zhHandoffBase
Set-MdHyperlink $wsZh "A3" $addedGuid $overviewRepo
Set-XlfHyperlink $wsZh "C3" $addedGuid $addedHash "zh-cn" $zhHandoffBase
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$overviewRepo/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A4").Value = ".localization-config"
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $deTime
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("D2").Value = $deTime

Set-MdHyperlink $wsDe "A2" $newGuid $overviewRepo
Set-XlfHyperlink $wsDe "C2" $newGuid $newHash "de-de" $deHandoffBase
Set-MdHyperlink $wsDe "A3" $addedGuid $overviewRepo
Set-XlfHyperlink $wsDe "C3" $addedGuid $addedHash "de-de" $deHandoffBase
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$overviewRepo/.localization-config", "", "", ".localization-config") | Out-Null
